$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D28").Value = "[24/02/14] Reward hacking이란 무엇인가?"
$ws.Range("E28").Value = "https://ropiens.tistory.com/242"

$ws.Range("D36").Value = "Controllable Diffusion Models"
$ws.Range("E36").Value = "http://dmqm.korea.ac.kr/activity/seminar/441"
